# Re-instate running without sheet behaviour
# Insert a new parameter row ("tb_multiplier_child_infectiousness") into the
# "constants" sheet just above the existing "tb_prop_early_progression" row
# (new row 6), so that the data-processing module always has a value to
# read for this parameter even when no country sheet is available.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Insert a new row before the current row 6 - this shifts every
# subsequent row (and any formulas referencing them) down by one,
# while inheriting the formatting of the row above it.
$ws.Rows.Item(6).Insert()

# Populate the new row with the new parameter.
$ws.Cells.Item(6, 1).Value2 = "tb_multiplier_child_infectiousness"
$ws.Cells.Item(6, 2).Value2 = 1
$ws.Cells.Item(6, 5).Value2 = "Note this is required for parameter loops, even though it is irrelevant epidemiologically"

# Move the selection/cursor to reflect where the author ended up editing.
[void]$ws.Range("E7").Select()
